$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing text (append new clauses describing popupData / popupHtml) ---
$ws.Range("B3").Value = 'В объектах в properties хранятся данные для фильтрации: .{ ownerLayer : "", ownerElement: "" ,  year: "", popupData : {}}'
$ws.Range("B4").Value = 'В слое в metadata хранится: {parentLayers: [], childLayers: [], elements: [], popupHtml: ""}'

# --- New cells documenting popupData / popupHtml ---
$ws.Range("H3").Value = "popupData - данные для отображения всплывающего окна при нажатии"
$ws.Range("F4").Value = "popupHtml: html модель всплывающего окна, с кнокаутовским датабинденгом"

# --- Column widths: widen column F slightly, add a custom width for the new column H ---
$ws.Columns.Item(6).ColumnWidth = 56.9
$ws.Columns.Item(8).ColumnWidth = 45.26

# --- Sheet view: scroll so column B is the leftmost visible column, move the selection ---
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C8").Select()
